$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a text value, forcing text type even when the
# string looks like a number (e.g. "1.002"), and restore the original
# "General" number format afterwards so no visible formatting changes.
function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

# --- Apply updated values (price & % volume changes; row 45/46 coin swap) ---
# Row 2
$ws.Range('D2').Value = '27.113.46'
$ws.Range('E2').Value = '  -2.71%  '

# Row 3
$ws.Range('D3').Value = '1.717.20'
$ws.Range('E3').Value = '  -2.94%  '

# Row 4
Set-TextValue $ws 'D4' '1.002'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
Set-TextValue $ws 'D5' '309.25'
$ws.Range('E5').Value = '  -5.68%  '

# Row 7
Set-TextValue $ws 'D7' '0.4716'
$ws.Range('E7').Value = '  +5.12%  '

# Row 8
Set-TextValue $ws 'D8' '0.3426'

# Row 9
Set-TextValue $ws 'D9' '42.11'
$ws.Range('E9').Value = '  +0.31%  '

# Row 10
Set-TextValue $ws 'D10' '0.07264'
$ws.Range('E10').Value = '  -2.11%  '

# Row 11
Set-TextValue $ws 'D11' '1.044'
$ws.Range('E11').Value = '  -4.89%  '

# Row 12
Set-TextValue $ws 'D12' '1.001'
$ws.Range('E12').Value = '  -0.05%  '

# Row 13
Set-TextValue $ws 'D13' '19.88'
$ws.Range('E13').Value = '  -5.03%  '

# Row 14
Set-TextValue $ws 'D14' '5.872'
$ws.Range('E14').Value = '  -2.52%  '

# Row 15
$ws.Range('D15').Value = '1.719.62'
$ws.Range('E15').Value = '  -2.99%  '

# Row 16
Set-TextValue $ws 'D16' '6.879'
$ws.Range('E16').Value = '  -5.04%  '

# Row 17
Set-TextValue $ws 'D17' '89.32'
$ws.Range('E17').Value = '  -4.17%  '

# Row 18
Set-TextValue $ws 'D18' '0.00001040'
$ws.Range('E18').Value = '  -1.88%  '

# Row 19
Set-TextValue $ws 'D19' '0.06350'
$ws.Range('E19').Value = '  -1.34%  '

# Row 20
Set-TextValue $ws 'D20' '1.000'
$ws.Range('E20').Value = '  -0.02%  '

# Row 21
$ws.Range('E21').Value = '  -3.45%  '

# Row 22
Set-TextValue $ws 'D22' '5.611'
$ws.Range('E22').Value = '  -2.88%  '

# Row 23
$ws.Range('D23').Value = '27.150.63'
$ws.Range('E23').Value = '  -2.69%  '

# Row 24
Set-TextValue $ws 'D24' '10.87'
$ws.Range('E24').Value = '  -3.67%  '

# Row 25
Set-TextValue $ws 'D25' '2.107'
$ws.Range('E25').Value = '  +0.01%  '

# Row 26
Set-TextValue $ws 'D26' '156.08'
$ws.Range('E26').Value = '  -4.09%  '

# Row 27
Set-TextValue $ws 'D27' '19.49'
$ws.Range('E27').Value = '  -4.30%  '

# Row 28
$ws.Range('D28').Value = '1.913.34'
$ws.Range('E28').Value = '  -3.22%  '

# Row 29
Set-TextValue $ws 'D29' '2.082'
$ws.Range('E29').Value = '  -3.47%  '

# Row 30
Set-TextValue $ws 'D30' '119.46'
$ws.Range('E30').Value = '  -4.16%  '

# Row 31
Set-TextValue $ws 'D31' '1.013'
$ws.Range('E31').Value = '  -8.53%  '

# Row 32
Set-TextValue $ws 'D32' '0.09168'
$ws.Range('E32').Value = '  -0.23%  '

# Row 33
Set-TextValue $ws 'D33' '3.593'
$ws.Range('E33').Value = '  -1.77%  '

# Row 34
Set-TextValue $ws 'D34' '5.319'
$ws.Range('E34').Value = '  -5.18%  '

# Row 35
Set-TextValue $ws 'D35' '0.02202'
$ws.Range('E35').Value = '  -3.81%  '

# Row 36
Set-TextValue $ws 'D36' '0.05814'
$ws.Range('E36').Value = '  -4.60%  '

# Row 37
Set-TextValue $ws 'D37' '11.01'
$ws.Range('E37').Value = '  -7.16%  '

# Row 38
Set-TextValue $ws 'D38' '0.1995'
$ws.Range('E38').Value = '  -5.02%  '

# Row 39
Set-TextValue $ws 'D39' '4.734'
$ws.Range('E39').Value = '  -4.49%  '

# Row 40
Set-TextValue $ws 'D40' '1.393'
$ws.Range('E40').Value = '  -0.45%  '

# Row 41
Set-TextValue $ws 'D41' '0.5897'
$ws.Range('E41').Value = '  -6.79%  '

# Row 42
Set-TextValue $ws 'D42' '1.116'
$ws.Range('E42').Value = '  -5.65%  '

# Row 43
Set-TextValue $ws 'D43' '7.477'
$ws.Range('E43').Value = '  -5.22%  '

# Row 44
Set-TextValue $ws 'D44' '12.62'
$ws.Range('E44').Value = '  -5.11%  '

# Row 45
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D45' '3.569'
$ws.Range('E45').Value = '  -4.64%  '

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws 'D46' '0.5653'
$ws.Range('E46').Value = '  -4.21%  '

# Row 47
Set-TextValue $ws 'D47' '117.48'
$ws.Range('E47').Value = '  -3.99%  '

# Row 48
Set-TextValue $ws 'D48' '1.843'
$ws.Range('E48').Value = '  -5.77%  '

# Row 49
Set-TextValue $ws 'D49' '0.06650'
$ws.Range('E49').Value = '  -3.68%  '

# Row 50
$ws.Range('E50').Value = '  -4.31%  '
